# Fix field validation and UI guide - update Input sheet to standard template format
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Start clean: remove all existing content/formatting in the old A1:Q3 range
$ws.Range("A1:Q3").Clear()

# ---- Header row (row 1), default (unstyled) formatting ----
$headers = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# ---- Data rows ----
# Columns A and B hold date-looking text ("2025-08-26" etc.). Force Text format
# while writing so the COM layer doesn't auto-convert them into date serials,
# then reset style back to Normal so no stray style id remains on the cells.
$dateRange = $ws.Range("A2:B3")
$dateRange.NumberFormat = "@"

$row2 = @("2025-08-26","2025-10-11","신호수","신호수@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","4. 장비비","2) 신호수 외","기타","5월","KS규격-1",1,155000,170500,"박수진 ")
$row3 = @("2025-09-03","2025-09-18","신호수","신호수@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","4. 장비비","2) 신호수 외","기타","5월","KS규격-2",1,155000,170500,"장승훈 ")

for ($col = 1; $col -le $row2.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $row2[$col - 1]
}
for ($col = 1; $col -le $row3.Length; $col++) {
    $ws.Cells.Item(3, $col).Value = $row3[$col - 1]
}

# Reset formatting on the whole populated range to the default/Normal style so
# no header bold/border/center styling (and no stray text-format style) remains.
$ws.Range("A1:P3").Style = "Normal"

Write-Host ("Final UsedRange: " + $ws.UsedRange.Address())
